$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.200.68"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "1.439.21"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9155"
$ws.Range("E5").Value = "  -8.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "275.49"
$ws.Range("E6").Value = "  +1.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3616"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3072"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "38.91"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.023"
$ws.Range("E10").Value = "  +2.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06479"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9991"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.338"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.47"
$ws.Range("E14").Value = "  +2.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.047"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001009"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "1.436.81"
$ws.Range("E17").Value = "  +3.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9353"
$ws.Range("E18").Value = "  -6.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05618"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.49"
$ws.Range("E20").Value = "  -4.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.381"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.26"
$ws.Range("E22").Value = "  -3.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.83"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.239"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").Value = "20.202.23"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.08"
$ws.Range("E26").Value = "  +2.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.130"
$ws.Range("E27").Value = "  -2.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.86"
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("D29").Value = "1.588.12"
$ws.Range("E29").Value = "  +2.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "109.69"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.827"
$ws.Range("E31").Value = "  -5.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8077"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.831"
$ws.Range("E33").Value = "  -7.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07631"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.472"
$ws.Range("E35").Value = "  +4.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05840"
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.649"
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.127"
$ws.Range("E38").Value = "  +5.49%  "
$ws.Range("E39").Value = "  -2.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.18"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1848"
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9279"
$ws.Range("E42").Value = "  -7.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.166"
$ws.Range("E43").Value = "  -13.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5209"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.487"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "116.88"
$ws.Range("E47").Value = "  +5.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5087"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.730"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06348"
$ws.Range("E50").Value = "  +3.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9874"
$ws.Range("E51").Value = "  -1.04%  "
